# Update cryptos list (prices and volume %) per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.824.85'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').Value = '2.342.21'
$ws.Range('E3').Value = '  +1.20%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'307.08"
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('D6').Value = "'101.45"
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').Value = "'0.509"
$ws.Range('E7').Value = '  -4.89%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = "'0.514"
$ws.Range('E9').Value = '  -2.00%  '
$ws.Range('D10').Value = "'35.10"
$ws.Range('E10').Value = '  -2.76%  '
$ws.Range('E11').Value = '  +0.70%  '
$ws.Range('D12').Value = "'0.0796"
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('D14').Value = "'6.83"
$ws.Range('E14').Value = '  -3.46%  '
$ws.Range('D15').Value = '2.716.93'
$ws.Range('E15').Value = '  +1.70%  '
$ws.Range('D16').Value = "'15.42"
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('D17').Value = '2.353.84'
$ws.Range('E17').Value = '  +1.82%  '
$ws.Range('D18').Value = "'0.798"
$ws.Range('E18').Value = '  -1.91%  '
$ws.Range('D19').Value = '42.804.41'
$ws.Range('E19').Value = '  -0.86%  '
$ws.Range('D20').Value = "'6.24"
$ws.Range('E20').Value = '  +1.47%  '
$ws.Range('D21').Value = "'11.72"
$ws.Range('E21').Value = '  -6.97%  '
$ws.Range('D22').Value = '0.0₃0904'
$ws.Range('E22').Value = '  -1.85%  '
$ws.Range('D23').Value = "'67.41"
$ws.Range('E23').Value = '  -1.51%  '
$ws.Range('D24').Value = "'237.00"
$ws.Range('E24').Value = '  -1.85%  '
$ws.Range('D25').Value = "'2.00"
$ws.Range('E25').Value = '  -1.65%  '
$ws.Range('D26').Value = "'2.57"
$ws.Range('E26').Value = '  -2.56%  '
$ws.Range('D27').Value = "'0.999"
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').Value = "'25.30"
$ws.Range('E28').Value = '  +2.35%  '
$ws.Range('D29').Value = "'3.84"
$ws.Range('E29').Value = '  -3.92%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = "'2.19"
$ws.Range('E30').Value = '  +3.70%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = "'35.36"
$ws.Range('E31').Value = '  -5.96%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = "'9.33"
$ws.Range('E32').Value = '  -3.45%  '
$ws.Range('D33').Value = "'160.28"
$ws.Range('E33').Value = '  -4.38%  '
$ws.Range('D34').Value = "'1.00"
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').Value = "'5.15"
$ws.Range('E35').Value = '  -3.71%  '
$ws.Range('D36').Value = "'17.82"
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('E37').Value = '  +3.93%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = "'0.0727"
$ws.Range('E38').Value = '  -2.46%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = "'4.57"
$ws.Range('E39').Value = '  +6.63%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = "'3.01"
$ws.Range('E40').Value = '  -4.56%  '
$ws.Range('D41').Value = "'1.88"
$ws.Range('E41').Value = '  +1.16%  '
$ws.Range('D42').Value = "'0.103"
$ws.Range('E42').Value = '  -3.48%  '
$ws.Range('E43').Value = '  -3.09%  '
$ws.Range('D44').Value = "'2.52"
$ws.Range('E44').Value = '  +9.22%  '
$ws.Range('D45').Value = '2.021.47'
$ws.Range('E45').Value = '  +2.29%  '
$ws.Range('D46').Value = "'19.20"
$ws.Range('E46').Value = '  -3.68%  '
$ws.Range('D47').Value = "'0.0285"
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('D48').Value = "'10.51"
$ws.Range('E48').Value = '  +7.15%  '
$ws.Range('D49').Value = "'3.00"
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('D50').Value = "'57.00"
$ws.Range('E50').Value = '  +2.23%  '
$ws.Range('D51').Value = "'2.92"
$ws.Range('E51').Value = '  -2.22%  '
